$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 13072
$ws1.Range("F5").Value  = 90
$ws1.Range("F8").Value  = 28
$ws1.Range("F10").Value = 13037
$ws1.Range("F11").Value = 297
$ws1.Range("F12").Value = 551
$ws1.Range("F13").Value = 8735
$ws1.Range("F14").Value = 7769
$ws1.Range("F15").Value = 211
$ws1.Range("F19").Value = 993
$ws1.Range("F20").Value = 11
$ws1.Range("F24").Value = 337

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 13072
$ws4.Range("F6").Value  = 90
$ws4.Range("F9").Value  = 28
$ws4.Range("F11").Value = 13037
$ws4.Range("F12").Value = 297
$ws4.Range("F13").Value = 551
$ws4.Range("F14").Value = 8735
$ws4.Range("F15").Value = 7769
$ws4.Range("F16").Value = 211
$ws4.Range("F20").Value = 993
$ws4.Range("F21").Value = 11
$ws4.Range("F27").Value = 337
